$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.504499188239663
$ws.Cells.Item(2, 3).Value = 0.1242063838438554
$ws.Cells.Item(2, 4).Value = 0.07855712732232689
$ws.Cells.Item(2, 5).Value = 0.1026203682630076
$ws.Cells.Item(2, 7).Value = 0.6269556958365783
$ws.Cells.Item(2, 8).Value = 0.7496865490056877
$ws.Cells.Item(2, 9).Value = 0.7343779279826492
$ws.Cells.Item(2, 11).Value = 0.3124922173064988
$ws.Cells.Item(2, 12).Value = 0.2017774848741567
$ws.Cells.Item(2, 14).Value = 1.531276483369242
$ws.Cells.Item(2, 15).Value = 2.739961954799384

$ws.Cells.Item(3, 2).Value = 0.4634117496288468
$ws.Cells.Item(3, 3).Value = 0.1229486988016717
$ws.Cells.Item(3, 4).Value = 0.07126311655875384
$ws.Cells.Item(3, 5).Value = 0.1021057143821515
$ws.Cells.Item(3, 7).Value = 0.6282197294803709
$ws.Cells.Item(3, 8).Value = 0.7538319354129612
$ws.Cells.Item(3, 9).Value = 0.7393243124759614
$ws.Cells.Item(3, 11).Value = 0.2755816036042233
$ws.Cells.Item(3, 12).Value = 0.1943617475487827
$ws.Cells.Item(3, 14).Value = 1.546579099026402
$ws.Cells.Item(3, 15).Value = 2.751065145700565

$ws.Cells.Item(4, 2).Value = 0.4382947548042466
$ws.Cells.Item(4, 3).Value = 0.1221711544935005
$ws.Cells.Item(4, 4).Value = 0.06681695194070869
$ws.Cells.Item(4, 5).Value = 0.1018408856788682
$ws.Cells.Item(4, 7).Value = 0.6293711392029024
$ws.Cells.Item(4, 8).Value = 0.7566732183091247
$ws.Cells.Item(4, 9).Value = 0.7427066019468604
$ws.Cells.Item(4, 11).Value = 0.2529002842930055
$ws.Cells.Item(4, 12).Value = 0.1899140448974634
$ws.Cells.Item(4, 14).Value = 1.556463899390764
$ws.Cells.Item(4, 15).Value = 2.759288743722749

$ws.Cells.Item(5, 2).Value = 0.4280879627122829
$ws.Cells.Item(5, 3).Value = 0.1218529826525057
$ws.Cells.Item(5, 4).Value = 0.0650132882266945
$ws.Cells.Item(5, 5).Value = 0.1017458577857511
$ws.Cells.Item(5, 7).Value = 0.6299347219129743
$ws.Cells.Item(5, 8).Value = 0.7579055768507388
$ws.Cells.Item(5, 9).Value = 0.7441717555247642
$ws.Cells.Item(5, 11).Value = 0.2436534772830043
$ws.Cells.Item(5, 12).Value = 0.1881282204366244
$ws.Cells.Item(5, 14).Value = 1.56061500309993
$ws.Cells.Item(5, 15).Value = 2.762993699602646

$ws.Cells.Item(6, 2).Value = 0.4263948801764172
$ws.Cells.Item(6, 3).Value = 0.1218000716028342
$ws.Cells.Item(6, 4).Value = 0.06471428674555568
$ws.Cells.Item(6, 5).Value = 0.1017308578421101
$ws.Cells.Item(6, 7).Value = 0.6300340044857577
$ws.Cells.Item(6, 8).Value = 0.7581147118276732
$ws.Cells.Item(6, 9).Value = 0.7444202895255962
$ws.Cells.Item(6, 11).Value = 0.2421178276398877
$ws.Cells.Item(6, 12).Value = 0.1878332981784752
$ws.Cells.Item(6, 14).Value = 1.561311718292082
$ws.Cells.Item(6, 15).Value = 2.763630276850492

$ws.Cells.Item(7, 2).Value = 0.4381569857465308
$ws.Cells.Item(7, 3).Value = 0.1221668688184465
$ws.Cells.Item(7, 4).Value = 0.06679259390902814
$ws.Cells.Item(7, 5).Value = 0.1018395518675455
$ws.Cells.Item(7, 7).Value = 0.62937835775422
$ws.Cells.Item(7, 8).Value = 0.7566895365208097
$ws.Cells.Item(7, 9).Value = 0.7427260098395365
$ws.Cells.Item(7, 11).Value = 0.2527755940216565
$ws.Cells.Item(7, 12).Value = 0.1898898526061004
$ws.Cells.Item(7, 14).Value = 1.556519384708125
$ws.Cells.Item(7, 15).Value = 2.759337277410651

$ws.Cells.Item(8, 2).Value = 0.490309644176989
$ws.Cells.Item(8, 3).Value = 0.1237738509685613
$ws.Cells.Item(8, 4).Value = 0.07603545247687293
$ws.Cells.Item(8, 5).Value = 0.1024323094693109
$ws.Cells.Item(8, 7).Value = 0.6273136450066019
$ws.Cells.Item(8, 8).Value = 0.7510544920850108
$ws.Cells.Item(8, 9).Value = 0.7360118224439773
$ws.Cells.Item(8, 11).Value = 0.2997695061804961
$ws.Cells.Item(8, 12).Value = 0.1991986801970569
$ws.Cells.Item(8, 14).Value = 1.536451361406886
$ws.Cells.Item(8, 15).Value = 2.74349857998935

$ws.Cells.Item(9, 2).Value = 0.59343439205918
$ws.Cells.Item(9, 3).Value = 0.1268821088722589
$ws.Cells.Item(9, 4).Value = 0.09441681251885825
$ws.Cells.Item(9, 5).Value = 0.1039998188238336
$ws.Cells.Item(9, 7).Value = 0.6262430279450939
$ws.Cells.Item(9, 8).Value = 0.7423494335597098
$ws.Cells.Item(9, 9).Value = 0.725582602947437
$ws.Cells.Item(9, 11).Value = 0.3917611990997045
$ws.Cells.Item(9, 12).Value = 0.2182878474459642
$ws.Cells.Item(9, 14).Value = 1.500975968291069
$ws.Cells.Item(9, 15).Value = 2.723590824067031

$ws.Cells.Item(10, 2).Value = 0.6696930570767563
$ws.Cells.Item(10, 3).Value = 0.1291386284523099
$ws.Cells.Item(10, 4).Value = 0.1080780064951909
$ws.Cells.Item(10, 5).Value = 0.1053974825715365
$ws.Cells.Item(10, 7).Value = 0.6272739071288242
$ws.Cells.Item(10, 8).Value = 0.7373793193466724
$ws.Cells.Item(10, 9).Value = 0.7195871027797551
$ws.Cells.Item(10, 11).Value = 0.4592275819496763
$ws.Cells.Item(10, 12).Value = 0.2328187697222859
$ws.Cells.Item(10, 14).Value = 1.477273382376405
$ws.Cells.Item(10, 15).Value = 2.715758482034261

$ws.Cells.Item(11, 2).Value = 0.7044867522394327
$ws.Cells.Item(11, 3).Value = 0.1301591142383742
$ws.Cells.Item(11, 4).Value = 0.1143269247939429
$ws.Cells.Item(11, 5).Value = 0.1060865487711062
$ws.Cells.Item(11, 7).Value = 0.6281379836534882
$ws.Cells.Item(11, 8).Value = 0.7354269512253353
$ws.Cells.Item(11, 9).Value = 0.717221169432726
$ws.Cells.Item(11, 11).Value = 0.4898897906615218
$ws.Cells.Item(11, 12).Value = 0.2395386842056269
$ws.Cells.Item(11, 14).Value = 1.467002557232522
$ws.Cells.Item(11, 15).Value = 2.713669816416143

$ws.Cells.Item(12, 2).Value = 0.7176764382623446
$ws.Cells.Item(12, 3).Value = 0.1305446611323759
$ws.Cells.Item(12, 4).Value = 0.1166981457889449
$ws.Cells.Item(12, 5).Value = 0.1063551179915514
$ws.Cells.Item(12, 7).Value = 0.6285220223834216
$ws.Cells.Item(12, 8).Value = 0.7347319378926045
$ws.Cells.Item(12, 9).Value = 0.7163771959842009
$ws.Cells.Item(12, 11).Value = 0.501496176494328
$ws.Cells.Item(12, 12).Value = 0.2420990402018077
$ws.Cells.Item(12, 14).Value = 1.463186815908291
$ws.Cells.Item(12, 15).Value = 2.713090785885839

$ws.Cells.Item(13, 2).Value = 0.714835190062729
$ws.Cells.Item(13, 3).Value = 0.130461666623944
$ws.Cells.Item(13, 4).Value = 0.1161872442122132
$ws.Cells.Item(13, 5).Value = 0.1062969376314911
$ws.Cells.Item(13, 7).Value = 0.628436784963597
$ws.Cells.Item(13, 8).Value = 0.7348796519513172
$ws.Cells.Item(13, 9).Value = 0.7165566504410421
$ws.Cells.Item(13, 11).Value = 0.4989967537013626
$ws.Cells.Item(13, 12).Value = 0.2415469265200301
$ws.Cells.Item(13, 14).Value = 1.464005330318908
$ws.Cells.Item(13, 15).Value = 2.713206067712974

$ws.Cells.Item(14, 2).Value = 0.7055715981454682
$ws.Cells.Item(14, 3).Value = 0.1301908513650218
$ws.Cells.Item(14, 4).Value = 0.1145219088424909
$ws.Cells.Item(14, 5).Value = 0.1061084912487189
$ws.Cells.Item(14, 7).Value = 0.6281684396856519
$ws.Cells.Item(14, 8).Value = 0.7353688844694943
$ws.Cells.Item(14, 9).Value = 0.7171506941159365
$ws.Cells.Item(14, 11).Value = 0.4908447532811806
$ws.Cells.Item(14, 12).Value = 0.2397490128018802
$ws.Cells.Item(14, 14).Value = 1.466687158720168
$ws.Cells.Item(14, 15).Value = 2.713617932984846

$ws.Cells.Item(15, 2).Value = 0.6998991896884377
$ws.Cells.Item(15, 3).Value = 0.1300248527290151
$ws.Cells.Item(15, 4).Value = 0.1135024782839764
$ws.Cells.Item(15, 5).Value = 0.1059940559546995
$ws.Cells.Item(15, 7).Value = 0.6280114720004093
$ws.Cells.Item(15, 8).Value = 0.7356743215600687
$ws.Cells.Item(15, 9).Value = 0.7175213283345485
$ws.Cells.Item(15, 11).Value = 0.4858507837224408
$ws.Cells.Item(15, 12).Value = 0.2386497762392423
$ws.Cells.Item(15, 14).Value = 1.468339441509835
$ws.Cells.Item(15, 15).Value = 2.713897805189191

$ws.Cells.Item(16, 2).Value = 0.6674212183644386
$ws.Cells.Item(16, 3).Value = 0.1290718142085723
$ws.Cells.Item(16, 4).Value = 0.1076703122754736
$ws.Cells.Item(16, 5).Value = 0.105353520176017
$ws.Cells.Item(16, 7).Value = 0.627225389764277
$ws.Cells.Item(16, 8).Value = 0.7375131147095715
$ws.Cells.Item(16, 9).Value = 0.7197489935135124
$ws.Cells.Item(16, 11).Value = 0.4572231065780556
$ws.Cells.Item(16, 12).Value = 0.2323818062325671
$ws.Cells.Item(16, 14).Value = 1.47795490214946
$ws.Cells.Item(16, 15).Value = 2.715924639395183

$ws.Cells.Item(17, 2).Value = 0.6475229254363342
$ws.Cells.Item(17, 3).Value = 0.1284855980810349
$ws.Cells.Item(17, 4).Value = 0.1041012346781258
$ws.Cells.Item(17, 5).Value = 0.1049741967110513
$ws.Cells.Item(17, 7).Value = 0.6268443643326407
$ws.Cells.Item(17, 8).Value = 0.7387201387044939
$ws.Cells.Item(17, 9).Value = 0.7212081543623903
$ws.Cells.Item(17, 11).Value = 0.4396531936692156
$ws.Cells.Item(17, 12).Value = 0.2285646361991951
$ws.Cells.Item(17, 14).Value = 1.483984729320367
$ws.Cells.Item(17, 15).Value = 2.717545569788882

$ws.Cells.Item(18, 2).Value = 0.6360877181699891
$ws.Cells.Item(18, 3).Value = 0.1281478572070966
$ws.Cells.Item(18, 4).Value = 0.1020516341010165
$ws.Cells.Item(18, 5).Value = 0.104761034609421
$ws.Cells.Item(18, 7).Value = 0.6266623988173166
$ws.Cells.Item(18, 8).Value = 0.739443434125576
$ws.Cells.Item(18, 9).Value = 0.7220814477618305
$ws.Cells.Item(18, 11).Value = 0.4295447905824403
$ws.Cells.Item(18, 12).Value = 0.2263794353120261
$ws.Cells.Item(18, 14).Value = 1.487501078715676
$ws.Cells.Item(18, 15).Value = 2.718616660890859

$ws.Cells.Item(19, 2).Value = 0.6322176542347222
$ws.Cells.Item(19, 3).Value = 0.1280334076609151
$ws.Cells.Item(19, 4).Value = 0.1013582323230935
$ws.Cells.Item(19, 5).Value = 0.1046897236028173
$ws.Cells.Item(19, 7).Value = 0.6266071754616149
$ws.Cells.Item(19, 8).Value = 0.7396933203787626
$ws.Cells.Item(19, 9).Value = 0.7223829739944208
$ws.Cells.Item(19, 11).Value = 0.4261218219467651
$ws.Cells.Item(19, 12).Value = 0.2256413420546437
$ws.Cells.Item(19, 14).Value = 1.488699923535517
$ws.Cells.Item(19, 15).Value = 2.7190031517205

$ws.Cells.Item(20, 2).Value = 0.6496401277273094
$ws.Cells.Item(20, 3).Value = 0.1285480604223679
$ws.Cells.Item(20, 4).Value = 0.1044808342834074
$ws.Cells.Item(20, 5).Value = 0.1050140575114789
$ws.Cells.Item(20, 7).Value = 0.6268810759913492
$ws.Cells.Item(20, 8).Value = 0.7385886431377742
$ws.Cells.Item(20, 9).Value = 0.7210493031379848
$ws.Cells.Item(20, 11).Value = 0.4415238200224962
$ws.Cells.Item(20, 12).Value = 0.2289699116798829
$ws.Cells.Item(20, 14).Value = 1.483337860133354
$ws.Cells.Item(20, 15).Value = 2.717358657039597

$ws.Cells.Item(21, 2).Value = 0.7082921640235327
$ws.Cells.Item(21, 3).Value = 0.1302704206850365
$ws.Cells.Item(21, 4).Value = 0.1150109258768737
$ws.Cells.Item(21, 5).Value = 0.1061636355065652
$ws.Cells.Item(21, 7).Value = 0.6282457167006612
$ws.Cells.Item(21, 8).Value = 0.7352239830792513
$ws.Cells.Item(21, 9).Value = 0.7169747992424647
$ws.Cells.Item(21, 11).Value = 0.4932393266203405
$ws.Cells.Item(21, 12).Value = 0.2402766792229016
$ws.Cells.Item(21, 14).Value = 1.46589744305334
$ws.Cells.Item(21, 15).Value = 2.713491208032764

$ws.Cells.Item(22, 2).Value = 0.7467063024987795
$ws.Cells.Item(22, 3).Value = 0.1313908903516889
$ws.Cells.Item(22, 4).Value = 0.1219214348188729
$ws.Cells.Item(22, 5).Value = 0.1069594413429655
$ws.Cells.Item(22, 7).Value = 0.6294688525324972
$ws.Cells.Item(22, 8).Value = 0.7332832032434595
$ws.Cells.Item(22, 9).Value = 0.7146146938152711
$ws.Cells.Item(22, 11).Value = 0.5270105598632142
$ws.Cells.Item(22, 12).Value = 0.2477575927592142
$ws.Cells.Item(22, 14).Value = 1.454928118840062
$ws.Cells.Item(22, 15).Value = 2.712198684371543

$ws.Cells.Item(23, 2).Value = 0.7261967407629868
$ws.Cells.Item(23, 3).Value = 0.1307933575467359
$ws.Cells.Item(23, 4).Value = 0.1182305770043968
$ws.Cells.Item(23, 5).Value = 0.1065306419684653
$ws.Cells.Item(23, 7).Value = 0.6287857270409347
$ws.Cells.Item(23, 8).Value = 0.73429542813561
$ws.Cells.Item(23, 9).Value = 0.7158466258435752
$ws.Cells.Item(23, 11).Value = 0.508988968257853
$ws.Cells.Item(23, 12).Value = 0.2437565715403309
$ws.Cells.Item(23, 14).Value = 1.460743389970725
$ws.Cells.Item(23, 15).Value = 2.71277555277976

$ws.Cells.Item(24, 2).Value = 0.6486829262270248
$ws.Cells.Item(24, 3).Value = 0.1285198234317733
$ws.Cells.Item(24, 4).Value = 0.1043092101044039
$ws.Cells.Item(24, 5).Value = 0.1049960211292884
$ws.Cells.Item(24, 7).Value = 0.6268643631156294
$ws.Cells.Item(24, 8).Value = 0.7386480008200351
$ws.Cells.Item(24, 9).Value = 0.7211210126138745
$ws.Cells.Item(24, 11).Value = 0.4406781323100688
$ws.Cells.Item(24, 12).Value = 0.2287866575317423
$ws.Cells.Item(24, 14).Value = 1.483630154814824
$ws.Cells.Item(24, 15).Value = 2.717442726697755

$ws.Cells.Item(25, 2).Value = 0.5654478714872369
$ws.Cells.Item(25, 3).Value = 0.1260459411181856
$ws.Cells.Item(25, 4).Value = 0.08941671564478781
$ws.Cells.Item(25, 5).Value = 0.1035324939571503
$ws.Cells.Item(25, 7).Value = 0.6262136222709813
$ws.Cells.Item(25, 8).Value = 0.7444537465661938
$ws.Cells.Item(25, 9).Value = 0.7281111321167835
$ws.Cells.Item(25, 11).Value = 0.3668945911471155
$ws.Cells.Item(25, 12).Value = 0.2130346477127745
$ws.Cells.Item(25, 14).Value = 1.510158152229959
$ws.Cells.Item(25, 15).Value = 2.713897805189191
